$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -19
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -5
$ws.Range("F8").Value = -4
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = -1
